# Apply red (FF0000) font color to several bullet paragraphs in the
# "design breakdown" (numId=4) list, and split the RJ45 paragraph's
# first run so only part of it becomes red, and merge/re-color the
# "controller that could be used ..." paragraph (including merging the
# "or  " / "PI" runs that Word's Find&Replace naturally coalesces).

$d = $word.ActiveDocument

# wdColorRed == 255 (0x0000FF in BGR == pure red in the OOXML <w:color w:val="FF0000"/> sense)
$wdColorRed = 255

# --- Paragraphs that become entirely red -------------------------------
# 17: "1x USB 3.2 Gen2 type C - from the hub controller"
# 18: "2x USB 3.2 gen1 type A - from the hub controller"
# 20: "1x type C connect to Ipad - via a basic type C with low-medium speed
#      capabilities (not high speed) this will be directly from the DP and
#      DM of the USB C which will be free."
# 27: "M.2 controller for NVME and SATA - We need a M.2 to USB
#      controller/interfacing IC - (yet to be establish the appropriate
#      controller)."
foreach ($idx in 17, 18, 20, 27) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.Color = $wdColorRed
}

# --- RJ45 paragraph: only the first run (up to "...Recommended") is red,
#     the following ": " stays black. -----------------------------------
$rj45 = $d.Paragraphs.Item(26)
$rj45Start = $rj45.Range.Start
$rj45End = $rj45.Range.End

$findRng = $d.Range($rj45Start, $rj45End)
$findRng.Find.ClearFormatting()
$findRng.Find.Execute("Recommended: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# findRng now covers "Recommended: " - split it 2 chars before the end
# (i.e. right before the ": ") so the colon+space is excluded.
$splitPoint = $findRng.End - 2
$redPart = $d.Range($rj45Start, $splitPoint)
$redPart.Font.Color = $wdColorRed

# --- "the controller that could be used ..." paragraph -----------------
# First merge the "or  " + "PI" runs (separated by a proofErr element)
# into a single "or  PI" run via a find & replace - Word coalesces the
# underlying runs that now share identical formatting.
$controllerPara = $d.Paragraphs.Item(28)
$mergeRng = $controllerPara.Range
$mergeRng.Find.ClearFormatting()
$mergeRng.Find.Execute("or  PI", $true, $false, $false, $false, $false, $true, 1, $false, "or  PI", 2)

# Now color the whole paragraph red (bold runs keep their bold).
$controllerPara2 = $d.Paragraphs.Item(28)
$controllerPara2.Range.Font.Color = $wdColorRed
